$d = $word.ActiveDocument

function Set-ParaText($para, [string]$newText) {
    $r = $para.Range
    $r2 = $d.Range($r.Start, $r.End - 1)
    $r2.Text = $newText
}

# ---------------------------------------------------------------------------
# 1. "Base build" feature list: shift each item's text into the previous
#    paragraph and give the last item the new trailing entry. (Food storage
#    drops off the top, "Larger buildings (4 tiles)" is appended at the end.)
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p4 = $d.Paragraphs(4)
$p5 = $d.Paragraphs(5)
$p6 = $d.Paragraphs(6)
$p7 = $d.Paragraphs(7)
$p8 = $d.Paragraphs(8)
$p9 = $d.Paragraphs(9)

$t4 = $p4.Range.Text.TrimEnd([char]13)
$t5 = $p5.Range.Text.TrimEnd([char]13)
$t6 = $p6.Range.Text.TrimEnd([char]13)
$t7 = $p7.Range.Text.TrimEnd([char]13)
$t8 = $p8.Range.Text.TrimEnd([char]13)
$t9 = $p9.Range.Text.TrimEnd([char]13)

Set-ParaText $p3 $t4
Set-ParaText $p4 $t5
Set-ParaText $p5 $t6
Set-ParaText $p6 $t7
Set-ParaText $p7 $t8
Set-ParaText $p8 $t9
Set-ParaText $p9 "Larger buildings (4 tiles)"

# ---------------------------------------------------------------------------
# 2. Bugs list: add a new bullet after "Hitting creatures from far away".
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Hitting creatures from far away")
$bugsPara = $d.Paragraphs(1)
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13) -eq "Hitting creatures from far away") {
        $bugsPara = $para
        break
    }
}
$bugsPara.Range.InsertParagraphAfter()

$newBugPara = $bugsPara.Next()
Set-ParaText $newBugPara "Directing follower with full inventory to gather resource sends them to the storage, then sets to idle (FIX LAST TARGET THING)"

# ---------------------------------------------------------------------------
# 3. "To do" list: rename + renumber "Implement food resource on HUD".
# ---------------------------------------------------------------------------
$todoPara = $d.Paragraphs(1)
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13) -eq "Implement food resource on HUD") {
        $todoPara = $para
        break
    }
}
Set-ParaText $todoPara "Re-order building menu (resource storage all together)"
$todoPara.Range.ListFormat.RemoveNumbers()
$todoPara.Range.ListFormat.ApplyNumberDefault()

# New trailing empty paragraph, indented, no list.
$todoPara.Range.InsertParagraphAfter()
$trailingPara = $todoPara.Next()
$trailingPara.Range.ListFormat.RemoveNumbers()
$trailingPara.Style = "Normal"
$trailingPara.Range.ParagraphFormat.LeftIndent = 18
